$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct rows whose match data was placed on the wrong line (re-paired by id) ---
# Row 9
$ws.Range("B9").Value2 = 6865285
$ws.Range("C9").Value2 = 'Bosnia Herzegovina Premier Liga'
$ws.Range("D9").Value2 = 'Bosnia  Herzegovina Premier Liga'
$ws.Range("E9").Value2 = 45150.5
$ws.Range("F9").Value2 = 'NK Igman Konjic'
$ws.Range("G9").Value2 = 'Sloga'
$ws.Range("H9").Value2 = 1
$ws.Range("I9").Value2 = 0
$ws.Range("J9").Value2 = 'H'
$ws.Range("K9").Value2 = 2
$ws.Range("L9").Value2 = 3.4
$ws.Range("M9").Value2 = 3.2
$ws.Range("N9").Value2 = 1.909
$ws.Range("O9").Value2 = 3.5
$ws.Range("P9").Value2 = 3.4
$ws.Range("Q9").Value2 = -0.5
$ws.Range("R9").Value2 = 1.95
$ws.Range("S9").Value2 = 1.85
$ws.Range("T9").Value2 = 2.5
$ws.Range("U9").Value2 = 1.85
$ws.Range("V9").Value2 = 1.95
$ws.Range("W9").Value2 = 0.909
$ws.Range("X9").Value2 = -1
$ws.Range("Y9").Value2 = -1
$ws.Range("Z9").Value2 = 0.95
$ws.Range("AA9").Value2 = -1
$ws.Range("AB9").Value2 = -1
$ws.Range("AC9").Value2 = 0.95

# Row 10
$ws.Range("B10").Value2 = 6865281
$ws.Range("C10").Value2 = 'Bosnia Herzegovina Premier Liga'
$ws.Range("D10").Value2 = 'Bosnia  Herzegovina Premier Liga'
$ws.Range("E10").Value2 = 45150.5
$ws.Range("F10").Value2 = 'GOSK Gabela'
$ws.Range("G10").Value2 = 'Zvijezda 09'
$ws.Range("H10").Value2 = 2
$ws.Range("I10").Value2 = 0
$ws.Range("J10").Value2 = 'H'
$ws.Range("K10").Value2 = 1.75
$ws.Range("L10").Value2 = 4
$ws.Range("M10").Value2 = 3.5
$ws.Range("N10").Value2 = 1.75
$ws.Range("O10").Value2 = 4
$ws.Range("P10").Value2 = 3.4
$ws.Range("Q10").Value2 = -0.5
$ws.Range("R10").Value2 = 1.8
$ws.Range("S10").Value2 = 2
$ws.Range("T10").Value2 = 2.5
$ws.Range("U10").Value2 = 1.85
$ws.Range("V10").Value2 = 1.95
$ws.Range("W10").Value2 = 0.75
$ws.Range("X10").Value2 = -1
$ws.Range("Y10").Value2 = -1
$ws.Range("Z10").Value2 = 0.8
$ws.Range("AA10").Value2 = -1
$ws.Range("AB10").Value2 = -1
$ws.Range("AC10").Value2 = 0.95

# Row 29
$ws.Range("B29").Value2 = 6865295
$ws.Range("C29").Value2 = 'Bosnia Herzegovina Premier Liga'
$ws.Range("D29").Value2 = 'Bosnia  Herzegovina Premier Liga'
$ws.Range("E29").Value2 = 45172.61458333334
$ws.Range("F29").Value2 = 'FK Tuzla City'
$ws.Range("G29").Value2 = 'NK Igman Konjic'
$ws.Range("H29").Value2 = 3
$ws.Range("I29").Value2 = 1
$ws.Range("J29").Value2 = 'H'
$ws.Range("K29").Value2 = 1.8
$ws.Range("L29").Value2 = 3.4
$ws.Range("M29").Value2 = 3.8
$ws.Range("N29").Value2 = 1.615
$ws.Range("O29").Value2 = 3.5
$ws.Range("P29").Value2 = 4.5
$ws.Range("Q29").Value2 = -0.75
$ws.Range("R29").Value2 = 1.85
$ws.Range("S29").Value2 = 1.95
$ws.Range("T29").Value2 = 2.75
$ws.Range("U29").Value2 = 2
$ws.Range("V29").Value2 = 1.8
$ws.Range("W29").Value2 = 0.615
$ws.Range("X29").Value2 = -1
$ws.Range("Y29").Value2 = -1
$ws.Range("Z29").Value2 = 0.8500000000000001
$ws.Range("AA29").Value2 = -1
$ws.Range("AB29").Value2 = 1
$ws.Range("AC29").Value2 = -1

# Row 30
$ws.Range("B30").Value2 = 6865296
$ws.Range("C30").Value2 = 'Bosnia Herzegovina Premier Liga'
$ws.Range("D30").Value2 = 'Bosnia  Herzegovina Premier Liga'
$ws.Range("E30").Value2 = 45172.61458333334
$ws.Range("F30").Value2 = 'Velez Mostar'
$ws.Range("G30").Value2 = 'Zeljeznicar'
$ws.Range("H30").Value2 = 1
$ws.Range("I30").Value2 = 0
$ws.Range("J30").Value2 = 'H'
$ws.Range("K30").Value2 = 1.909
$ws.Range("L30").Value2 = 3.2
$ws.Range("M30").Value2 = 3.6
$ws.Range("N30").Value2 = 1.95
$ws.Range("O30").Value2 = 3.2
$ws.Range("P30").Value2 = 3.4
$ws.Range("Q30").Value2 = -0.5
$ws.Range("R30").Value2 = 2.025
$ws.Range("S30").Value2 = 1.775
$ws.Range("T30").Value2 = 2.25
$ws.Range("U30").Value2 = 1.9
$ws.Range("V30").Value2 = 1.9
$ws.Range("W30").Value2 = 0.95
$ws.Range("X30").Value2 = -1
$ws.Range("Y30").Value2 = -1
$ws.Range("Z30").Value2 = 1.025
$ws.Range("AA30").Value2 = -1
$ws.Range("AB30").Value2 = -1
$ws.Range("AC30").Value2 = 0.8999999999999999

# Row 36
$ws.Range("B36").Value2 = 6865299
$ws.Range("C36").Value2 = 'Bosnia Herzegovina Premier Liga'
$ws.Range("D36").Value2 = 'Bosnia  Herzegovina Premier Liga'
$ws.Range("E36").Value2 = 45186.61458333334
$ws.Range("F36").Value2 = 'Siroki Brijeg'
$ws.Range("G36").Value2 = 'Zvijezda 09'
$ws.Range("H36").Value2 = 2
$ws.Range("I36").Value2 = 1
$ws.Range("J36").Value2 = 'H'
$ws.Range("K36").Value2 = 1.25
$ws.Range("L36").Value2 = 5.5
$ws.Range("M36").Value2 = 8
$ws.Range("N36").Value2 = 1.4
$ws.Range("O36").Value2 = 4.75
$ws.Range("P36").Value2 = 5.75
$ws.Range("Q36").Value2 = -1.25
$ws.Range("R36").Value2 = 1.9
$ws.Range("S36").Value2 = 1.9
$ws.Range("T36").Value2 = 2.75
$ws.Range("U36").Value2 = 1.85
$ws.Range("V36").Value2 = 1.95
$ws.Range("W36").Value2 = 0.3999999999999999
$ws.Range("X36").Value2 = -1
$ws.Range("Y36").Value2 = -1
$ws.Range("Z36").Value2 = -0.5
$ws.Range("AA36").Value2 = 0.45
$ws.Range("AB36").Value2 = 0.425
$ws.Range("AC36").Value2 = -0.5

# Row 37
$ws.Range("B37").Value2 = 6864629
$ws.Range("C37").Value2 = 'Bosnia Herzegovina Premier Liga'
$ws.Range("D37").Value2 = 'Bosnia  Herzegovina Premier Liga'
$ws.Range("E37").Value2 = 45186.61458333334
$ws.Range("F37").Value2 = 'Borac Banja Luka'
$ws.Range("G37").Value2 = 'NK Posusje'
$ws.Range("H37").Value2 = 1
$ws.Range("I37").Value2 = 0
$ws.Range("J37").Value2 = 'H'
$ws.Range("K37").Value2 = 1.363
$ws.Range("L37").Value2 = 4.5
$ws.Range("M37").Value2 = 6.5
$ws.Range("N37").Value2 = 1.363
$ws.Range("O37").Value2 = 4.2
$ws.Range("P37").Value2 = 6.5
$ws.Range("Q37").Value2 = -1.25
$ws.Range("R37").Value2 = 1.95
$ws.Range("S37").Value2 = 1.85
$ws.Range("T37").Value2 = 2.5
$ws.Range("U37").Value2 = 1.925
$ws.Range("V37").Value2 = 1.875
$ws.Range("W37").Value2 = 0.363
$ws.Range("X37").Value2 = -1
$ws.Range("Y37").Value2 = -1
$ws.Range("Z37").Value2 = -0.5
$ws.Range("AA37").Value2 = 0.425
$ws.Range("AB37").Value2 = -1
$ws.Range("AC37").Value2 = 0.875

# Row 49
$ws.Range("B49").Value2 = 6865311
$ws.Range("C49").Value2 = 'Bosnia Herzegovina Premier Liga'
$ws.Range("D49").Value2 = 'Bosnia  Herzegovina Premier Liga'
$ws.Range("E49").Value2 = 45200.41666666666
$ws.Range("F49").Value2 = 'Sloga'
$ws.Range("G49").Value2 = 'GOSK Gabela'
$ws.Range("H49").Value2 = 3
$ws.Range("I49").Value2 = 2
$ws.Range("J49").Value2 = 'H'
$ws.Range("K49").Value2 = 1.833
$ws.Range("L49").Value2 = 3.6
$ws.Range("M49").Value2 = 3.4
$ws.Range("N49").Value2 = 1.909
$ws.Range("O49").Value2 = 3.4
$ws.Range("P49").Value2 = 3.3
$ws.Range("Q49").Value2 = -0.5
$ws.Range("R49").Value2 = 1.925
$ws.Range("S49").Value2 = 1.875
$ws.Range("T49").Value2 = 2.25
$ws.Range("U49").Value2 = 1.825
$ws.Range("V49").Value2 = 1.975
$ws.Range("W49").Value2 = 0.909
$ws.Range("X49").Value2 = -1
$ws.Range("Y49").Value2 = -1
$ws.Range("Z49").Value2 = 0.925
$ws.Range("AA49").Value2 = -1
$ws.Range("AB49").Value2 = 0.825
$ws.Range("AC49").Value2 = -1

# Row 50
$ws.Range("B50").Value2 = 6865310
$ws.Range("C50").Value2 = 'Bosnia Herzegovina Premier Liga'
$ws.Range("D50").Value2 = 'Bosnia  Herzegovina Premier Liga'
$ws.Range("E50").Value2 = 45200.41666666666
$ws.Range("F50").Value2 = 'NK Igman Konjic'
$ws.Range("G50").Value2 = 'Zrinjski Mostar'
$ws.Range("H50").Value2 = 0
$ws.Range("I50").Value2 = 2
$ws.Range("J50").Value2 = 'A'
$ws.Range("K50").Value2 = 3.4
$ws.Range("L50").Value2 = 3.6
$ws.Range("M50").Value2 = 1.833
$ws.Range("N50").Value2 = 4.75
$ws.Range("O50").Value2 = 4.75
$ws.Range("P50").Value2 = 1.45
$ws.Range("Q50").Value2 = 1.25
$ws.Range("R50").Value2 = 1.775
$ws.Range("S50").Value2 = 2.025
$ws.Range("T50").Value2 = 2.75
$ws.Range("U50").Value2 = 1.85
$ws.Range("V50").Value2 = 1.95
$ws.Range("W50").Value2 = -1
$ws.Range("X50").Value2 = -1
$ws.Range("Y50").Value2 = 0.45
$ws.Range("Z50").Value2 = -1
$ws.Range("AA50").Value2 = 1.025
$ws.Range("AB50").Value2 = -1
$ws.Range("AC50").Value2 = 0.95

# Row 76
$ws.Range("B76").Value2 = 6865377
$ws.Range("C76").Value2 = 'Bosnia Herzegovina Premier Liga'
$ws.Range("D76").Value2 = 'Bosnia  Herzegovina Premier Liga'
$ws.Range("E76").Value2 = 45235.375
$ws.Range("F76").Value2 = 'Zrinjski Mostar'
$ws.Range("G76").Value2 = 'FK Tuzla City'
$ws.Range("H76").Value2 = 3
$ws.Range("I76").Value2 = 1
$ws.Range("J76").Value2 = 'H'
$ws.Range("K76").Value2 = 1.333
$ws.Range("L76").Value2 = 5
$ws.Range("M76").Value2 = 6
$ws.Range("N76").Value2 = 1.166
$ws.Range("O76").Value2 = 6.5
$ws.Range("P76").Value2 = 13
$ws.Range("Q76").Value2 = -2
$ws.Range("R76").Value2 = 1.9
$ws.Range("S76").Value2 = 1.9
$ws.Range("T76").Value2 = 3.25
$ws.Range("U76").Value2 = 1.95
$ws.Range("V76").Value2 = 1.85
$ws.Range("W76").Value2 = 0.1659999999999999
$ws.Range("X76").Value2 = -1
$ws.Range("Y76").Value2 = -1
$ws.Range("Z76").Value2 = 0
$ws.Range("AA76").Value2 = -0
$ws.Range("AB76").Value2 = 0.95
$ws.Range("AC76").Value2 = -1

# Row 77
$ws.Range("B77").Value2 = 6865328
$ws.Range("C77").Value2 = 'Bosnia Herzegovina Premier Liga'
$ws.Range("D77").Value2 = 'Bosnia  Herzegovina Premier Liga'
$ws.Range("E77").Value2 = 45235.375
$ws.Range("F77").Value2 = 'Siroki Brijeg'
$ws.Range("G77").Value2 = 'NK Posusje'
$ws.Range("H77").Value2 = 1
$ws.Range("I77").Value2 = 1
$ws.Range("J77").Value2 = 'D'
$ws.Range("K77").Value2 = 2
$ws.Range("L77").Value2 = 3
$ws.Range("M77").Value2 = 3.5
$ws.Range("N77").Value2 = 2.1
$ws.Range("O77").Value2 = 3
$ws.Range("P77").Value2 = 3.3
$ws.Range("Q77").Value2 = -0.25
$ws.Range("R77").Value2 = 1.825
$ws.Range("S77").Value2 = 1.975
$ws.Range("T77").Value2 = 2
$ws.Range("U77").Value2 = 1.825
$ws.Range("V77").Value2 = 1.975
$ws.Range("W77").Value2 = -1
$ws.Range("X77").Value2 = 2
$ws.Range("Y77").Value2 = -1
$ws.Range("Z77").Value2 = -0.5
$ws.Range("AA77").Value2 = 0.4875
$ws.Range("AB77").Value2 = 0
$ws.Range("AC77").Value2 = -0

# Row 111
$ws.Range("B111").Value2 = 6865352
$ws.Range("C111").Value2 = 'Bosnia Herzegovina Premier Liga'
$ws.Range("D111").Value2 = 'Bosnia  Herzegovina Premier Liga'
$ws.Range("E111").Value2 = 45339.375
$ws.Range("F111").Value2 = 'NK Posusje'
$ws.Range("G111").Value2 = 'Zvijezda 09'
$ws.Range("H111").Value2 = 2
$ws.Range("I111").Value2 = 0
$ws.Range("J111").Value2 = 'H'
$ws.Range("K111").Value2 = 1.615
$ws.Range("L111").Value2 = 3.5
$ws.Range("M111").Value2 = 4.75
$ws.Range("N111").Value2 = 1.5
$ws.Range("O111").Value2 = 3.6
$ws.Range("P111").Value2 = 5.75
$ws.Range("Q111").Value2 = -1
$ws.Range("R111").Value2 = 1.9
$ws.Range("S111").Value2 = 1.9
$ws.Range("T111").Value2 = 2.25
$ws.Range("U111").Value2 = 1.85
$ws.Range("V111").Value2 = 1.95
$ws.Range("W111").Value2 = 0.5
$ws.Range("X111").Value2 = -1
$ws.Range("Y111").Value2 = -1
$ws.Range("Z111").Value2 = 0.8999999999999999
$ws.Range("AA111").Value2 = -1
$ws.Range("AB111").Value2 = -0.5
$ws.Range("AC111").Value2 = 0.475

# Row 112
$ws.Range("B112").Value2 = 6865354
$ws.Range("C112").Value2 = 'Bosnia Herzegovina Premier Liga'
$ws.Range("D112").Value2 = 'Bosnia  Herzegovina Premier Liga'
$ws.Range("E112").Value2 = 45339.375
$ws.Range("F112").Value2 = 'NK Igman Konjic'
$ws.Range("G112").Value2 = 'GOSK Gabela'
$ws.Range("H112").Value2 = 1
$ws.Range("I112").Value2 = 2
$ws.Range("J112").Value2 = 'A'
$ws.Range("K112").Value2 = 1.8
$ws.Range("L112").Value2 = 3.25
$ws.Range("M112").Value2 = 4
$ws.Range("N112").Value2 = 2.25
$ws.Range("O112").Value2 = 3.1
$ws.Range("P112").Value2 = 2.9
$ws.Range("Q112").Value2 = -0.25
$ws.Range("R112").Value2 = 1.975
$ws.Range("S112").Value2 = 1.825
$ws.Range("T112").Value2 = 2.25
$ws.Range("U112").Value2 = 1.875
$ws.Range("V112").Value2 = 1.925
$ws.Range("W112").Value2 = -1
$ws.Range("X112").Value2 = -1
$ws.Range("Y112").Value2 = 1.9
$ws.Range("Z112").Value2 = -1
$ws.Range("AA112").Value2 = 0.825
$ws.Range("AB112").Value2 = 0.875
$ws.Range("AC112").Value2 = -1

# --- Append new match rows 140-143 ---
# Row 140
$ws.Range("A140").Value2 = 138
$ws.Range("B140").Value2 = 7952735
$ws.Range("C140").Value2 = 'Bosnia Herzegovina Premier Liga'
$ws.Range("D140").Value2 = 'Bosnia  Herzegovina Premier Liga'
$ws.Range("E140").Value2 = 45380.66666666666
$ws.Range("F140").Value2 = 'Zeljeznicar'
$ws.Range("G140").Value2 = 'Velez Mostar'
$ws.Range("K140").Value2 = 2.4
$ws.Range("L140").Value2 = 3.2
$ws.Range("M140").Value2 = 2.6
$ws.Range("N140").Value2 = 2.3
$ws.Range("O140").Value2 = 3.25
$ws.Range("P140").Value2 = 2.8
$ws.Range("Q140").Value2 = -0.25
$ws.Range("R140").Value2 = 2.05
$ws.Range("S140").Value2 = 1.75
$ws.Range("T140").Value2 = 1.75
$ws.Range("U140").Value2 = 1.75
$ws.Range("V140").Value2 = 2.05
$ws.Range("W140").Value2 = 0
$ws.Range("X140").Value2 = 0
$ws.Range("Y140").Value2 = 0
$ws.Range("Z140").Value2 = 0
$ws.Range("AA140").Value2 = 0
$ws.Range("A2").Copy()
$ws.Range("A140").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E140").PasteSpecial(-4122)

# Row 141
$ws.Range("A141").Value2 = 139
$ws.Range("B141").Value2 = 7952739
$ws.Range("C141").Value2 = 'Bosnia Herzegovina Premier Liga'
$ws.Range("D141").Value2 = 'Bosnia  Herzegovina Premier Liga'
$ws.Range("E141").Value2 = 45381.39583333334
$ws.Range("F141").Value2 = 'Zvijezda 09'
$ws.Range("G141").Value2 = 'Siroki Brijeg'
$ws.Range("K141").Value2 = 2.25
$ws.Range("L141").Value2 = 3.1
$ws.Range("M141").Value2 = 2.875
$ws.Range("N141").Value2 = 2.15
$ws.Range("O141").Value2 = 3
$ws.Range("P141").Value2 = 3.1
$ws.Range("Q141").Value2 = -0.25
$ws.Range("R141").Value2 = 1.95
$ws.Range("S141").Value2 = 1.85
$ws.Range("T141").Value2 = 2.25
$ws.Range("U141").Value2 = 2.05
$ws.Range("V141").Value2 = 1.75
$ws.Range("W141").Value2 = 0
$ws.Range("X141").Value2 = 0
$ws.Range("Y141").Value2 = 0
$ws.Range("Z141").Value2 = 0
$ws.Range("AA141").Value2 = 0
$ws.Range("A2").Copy()
$ws.Range("A141").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E141").PasteSpecial(-4122)

# Row 142
$ws.Range("A142").Value2 = 140
$ws.Range("B142").Value2 = 7952456
$ws.Range("C142").Value2 = 'Bosnia Herzegovina Premier Liga'
$ws.Range("D142").Value2 = 'Bosnia  Herzegovina Premier Liga'
$ws.Range("E142").Value2 = 45381.5
$ws.Range("F142").Value2 = 'Borac Banja Luka'
$ws.Range("G142").Value2 = 'Zrinjski Mostar'
$ws.Range("K142").Value2 = 2.2
$ws.Range("L142").Value2 = 3.2
$ws.Range("M142").Value2 = 2.875
$ws.Range("N142").Value2 = 2.45
$ws.Range("O142").Value2 = 3
$ws.Range("P142").Value2 = 2.7
$ws.Range("Q142").Value2 = 0
$ws.Range("R142").Value2 = 1.825
$ws.Range("S142").Value2 = 1.975
$ws.Range("T142").Value2 = 2
$ws.Range("U142").Value2 = 1.775
$ws.Range("V142").Value2 = 2.025
$ws.Range("W142").Value2 = 0
$ws.Range("X142").Value2 = 0
$ws.Range("Y142").Value2 = 0
$ws.Range("Z142").Value2 = 0
$ws.Range("AA142").Value2 = 0
$ws.Range("A2").Copy()
$ws.Range("A142").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E142").PasteSpecial(-4122)

# Row 143
$ws.Range("A143").Value2 = 141
$ws.Range("B143").Value2 = 7952737
$ws.Range("C143").Value2 = 'Bosnia Herzegovina Premier Liga'
$ws.Range("D143").Value2 = 'Bosnia  Herzegovina Premier Liga'
$ws.Range("E143").Value2 = 45381.69791666666
$ws.Range("F143").Value2 = 'GOSK Gabela'
$ws.Range("G143").Value2 = 'NK Posusje'
$ws.Range("K143").Value2 = 2.5
$ws.Range("L143").Value2 = 3.2
$ws.Range("M143").Value2 = 2.5
$ws.Range("N143").Value2 = 3
$ws.Range("O143").Value2 = 3.25
$ws.Range("P143").Value2 = 2.2
$ws.Range("Q143").Value2 = 0.25
$ws.Range("R143").Value2 = 1.85
$ws.Range("S143").Value2 = 1.95
$ws.Range("T143").Value2 = 2.25
$ws.Range("U143").Value2 = 1.85
$ws.Range("V143").Value2 = 1.95
$ws.Range("W143").Value2 = 0
$ws.Range("X143").Value2 = 0
$ws.Range("Y143").Value2 = 0
$ws.Range("Z143").Value2 = 0
$ws.Range("AA143").Value2 = 0
$ws.Range("A2").Copy()
$ws.Range("A143").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E143").PasteSpecial(-4122)

$excel.CutCopyMode = 0
